$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo "locacalizacion" -> "localizacion" in B1
$ws.Range("B1").Value = "localizacion"

# Update the selection to B1 (was A2)
$ws.Range("B1").Select()
